$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

for ($r = 2; $r -le 12; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $a.NumberFormat = "@"
    $a.Value = "23"
    $a.Style = "Normal"

    $m = $ws.Cells.Item($r, 13)
    $m.Value = "nan"
}
